$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# --- Update the query timestamps in column F of the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:20:17.404547"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:17.404554"
$dataSheet.Range("F4").Value = "2021-10-05 14:20:17.404558"
$dataSheet.Range("F5").Value = "2021-10-05 14:20:17.404560"
$dataSheet.Range("F6").Value = "2021-10-05 14:20:17.404563"
$dataSheet.Range("F7").Value = "2021-10-05 14:20:17.404566"
$dataSheet.Range("F8").Value = "2021-10-05 14:20:17.404568"
$dataSheet.Range("F9").Value = "2021-10-05 14:20:17.404571"

# --- Add the new "metadata" sheet right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold/centered/bordered style, matching "data"'s header style)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Copy the header cell style from "data"!F1 onto metadata!B1:G1
$dataSheet.Range("F1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122) | Out-Null

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Familial melanoma"
$metaSheet.Range("C2").Value = 522
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.10"
$metaSheet.Range("D2").Style = "Normal"
$metaSheet.Range("E2").Value = "2021-04-14T09:18:05.194303Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:20:17.401223"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/522/?format=json"

# A2 uses the same style as "data"!A2 (the index/id column style)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122) | Out-Null

# Keep "data" as the active sheet/tab, matching the original workbook state
$dataSheet.Activate()
$dataSheet.Range("A1").Select() | Out-Null
